# Update the cryptos list (price / volume refresh) as produced by the
# scheduled GitHub Actions run on Tue Sep 19 22:36:52 UTC 2023.
#
# Notes:
#  - Column D ("Price") holds values that are stored as *text* in the
#    workbook (they use "." as a thousands separator, e.g. "27.183.63",
#    so they must never be re-interpreted as numbers). For values that
#    Excel's automatic type detection would otherwise read as a genuine
#    number (e.g. "217.23", "20.07") we prefix the literal with a single
#    leading apostrophe so Excel keeps storing it as text instead of
#    silently converting it to a floating point number.
#  - Rows 38/39 and 50/51 are not simple value tweaks: the two coins in
#    each pair swap places (new coin/link/price/volume values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.183.63"
$ws.Range("E2").Value = "  +1.20%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.641.82"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'217.23"
$ws.Range("E5").Value = "  +0.35%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +1.05%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.93%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'20.07"
$ws.Range("E10").Value = "  +1.43%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.871.99"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.645.45"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.69%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +2.92%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'67.37"
$ws.Range("E16").Value = "  +1.62%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.161.83"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'218.85"
$ws.Range("E19").Value = "  -0.23%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'6.87"
$ws.Range("E21").Value = "  +3.74%  "

# Row 22 - Toncoin
$ws.Range("E22").Value = "  +6.73%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +0.72%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +0.36%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'147.73"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.54"
$ws.Range("E26").Value = "  +1.88%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.22%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.07%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.67%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.49%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.86%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.04%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.21%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.272.44"
$ws.Range("E35").Value = "  +2.15%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.99%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +1.91%  "

# Row 38 - was ARBITRUM, now ImmutableX (rows 38/39 swap places)
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.546"
$ws.Range("E38").Value = "  +1.15%  "

# Row 39 - was ImmutableX, now ARBITRUM
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.855"
$ws.Range("E39").Value = "  +2.84%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.03%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.810"
$ws.Range("E41").Value = "  +0.46%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  +7.04%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'5.31"
$ws.Range("E43").Value = "  -1.08%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.782.01"
$ws.Range("E44").Value = "  -0.09%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'61.77"
$ws.Range("E45").Value = "  +1.71%  "

# Row 46 - Quant
$ws.Range("D46").Value = "'91.86"
$ws.Range("E46").Value = "  +0.42%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +1.66%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  +0.80%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.17%  "

# Row 50 - was Algorand, now EnergySwap (rows 50/51 swap places)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.69"
$ws.Range("E50").Value = "  +1.19%  "

# Row 51 - was EnergySwap, now Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0976"
$ws.Range("E51").Value = "  +0.07%  "
